$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 593; this pushes existing rows 593-640 down to 594-641
$ws.Rows.Item(593).Insert()

# Populate the new row 593 with the new weekly price record
$ws.Cells.Item(593, 1).Value = 3
$ws.Cells.Item(593, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(593, 3).Value = "Coquimbo"
$ws.Cells.Item(593, 4).Value = 45223
$ws.Cells.Item(593, 5).Value = 5
$ws.Cells.Item(593, 6).Value = 100112009
$ws.Cells.Item(593, 7).Value = "Acelga"
$ws.Cells.Item(593, 8).Value = "Sin especificar"
$ws.Cells.Item(593, 9).Value = "Primera"
$ws.Cells.Item(593, 10).Value = 230
$ws.Cells.Item(593, 11).Value = 3000
$ws.Cells.Item(593, 12).Value = 3300
$ws.Cells.Item(593, 13).Value = 3104
$ws.Cells.Item(593, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(593, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(593, 16).Value = 517
$ws.Cells.Item(593, 17).Value = 6
$ws.Cells.Item(593, 18).Value = "Hortaliza"
